$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): table group names ---
$ws.Range("A1").Value = "user"
$ws.Range("C1").Value = "field"
$ws.Range("D1").Value = "field_images"
$ws.Range("E1").Value = "field_rounds"
$ws.Range("F1").Value = "fields_bookings"

# --- Row 3: id column for every table ---
$ws.Range("A3").Value = "id"
$ws.Range("C3").Value = "id"
$ws.Range("D3").Value = "id"
$ws.Range("E3").Value = "id"
$ws.Range("F3").Value = "id"

# --- Row 4 ---
$ws.Range("A4").Value = "username"
$ws.Range("C4").Value = "name"
$ws.Range("D4").Value = "image_name"
$ws.Range("E4").Value = "start_at"
$ws.Range("F4").Value = "user_id"

# --- Row 5 ---
$ws.Range("A5").Value = "firstname"
$ws.Range("C5").Value = "description"
$ws.Range("D5").Value = "field_id"
$ws.Range("E5").Value = "end_after"
$ws.Range("F5").Value = "field_round_id"

# --- Row 6 ---
$ws.Range("A6").Value = "lastname"
$ws.Range("C6").Value = "capacity"
$ws.Range("E6").Value = "price"

# --- Row 7 ---
$ws.Range("A7").Value = "address"
$ws.Range("C7").Value = "indoor"
$ws.Range("E7").Value = "field_id"

# --- Row 8 ---
$ws.Range("A8").Value = "email"
$ws.Range("C8").Value = "floor"
$ws.Range("E8").Value = "status *"

# --- Row 9 ---
$ws.Range("A9").Value = "phone"
$ws.Range("C9").Value = "address"

# --- Row 10 ---
$ws.Range("A10").Value = "password"
$ws.Range("C10").Value = "user_id"

# --- Row 11 ---
$ws.Range("A11").Value = "role *"

# --- Column widths (target stored widths: A=14.42578125, C=16, D=13.7109375, E=15.140625, F=16.140625;
#     inputs chosen as the closest achievable given the engine's internal width quantization) ---
$ws.Columns.Item(1).ColumnWidth = 13.6666666666667
$ws.Columns.Item(3).ColumnWidth = 15.1666666666667
$ws.Columns.Item(4).ColumnWidth = 12.8333333333333
$ws.Columns.Item(5).ColumnWidth = 14.3333333333333
$ws.Columns.Item(6).ColumnWidth = 15.3333333333333

# --- Selection ---
[void]$ws.Range("D6").Select()
